$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 - this shifts existing rows 2-11 down to 3-12,
# carrying their values, shared-string references, and cell styles with them
# (including the two hyperlink-styled Password cells, which keep style "3").
$ws.Rows(2).Insert()

# Populate the newly inserted row 2 with the new test case.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "cus01"
$ws.Range("C2").Value = "MatKhau123"
$ws.Range("D2").Value = "Account"

# New row should look like a normal data row (no special style), except the
# Password cell which uses the same "hyperlink-like" style as the other
# Password cells in column C.
$ws.Range("A2:D2").Style = "Normal"
$ws.Range("C2").Style = $ws.Range("C3").Style

# Renumber the "Stt" column for the rows that got pushed down (they kept
# their original numbers when Insert() shifted them).
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11

# The row insert does not move the worksheet's hyperlink anchors, so redo
# the two hyperlinks at their new (shifted-down) locations.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:MatKhau@123")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:MatKhau@123", [Type]::Missing, [Type]::Missing, "MatKhau@123")

# Match the saved cursor position and sheet dimension from the edited file.
$ws.Range("G10").Select()
